$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.446800000000003
$ws.Range("A8").Value = -22.43700000000002
$ws.Range("A10").Value = -22.09630000000001
$ws.Range("B11").Value = 6.569100000000001
$ws.Range("A12").Value = -21.42799999999999
$ws.Range("B12").Value = 4.793199999999996
$ws.Range("B15").Value = 5.945499999999994
$ws.Range("B17").Value = 4.040600000000005
$ws.Range("A18").Value = -22.36380000000001
$ws.Range("A25").Value = -21.91979999999998
$ws.Range("B26").Value = 5.317199999999998
$ws.Range("B27").Value = 6.791500000000005
$ws.Range("B28").Value = 6.4034
$ws.Range("B32").Value = 7.861800000000001
$ws.Range("A37").Value = -22.07460000000001
$ws.Range("B37").Value = 5.951500000000002
$ws.Range("B41").Value = 9.079700000000004
$ws.Range("B47").Value = 6.956900000000003
$ws.Range("B51").Value = 5.8952
$ws.Range("A55").Value = -22.09139999999999
$ws.Range("B65").Value = 5.841399999999998
$ws.Range("A68").Value = -21.46259999999999
$ws.Range("B73").Value = 8.937099999999997
$ws.Range("A77").Value = -20.22089999999999
$ws.Range("A78").Value = -19.48219999999998
$ws.Range("A79").Value = -20.58609999999999
$ws.Range("A80").Value = -19.5015
$ws.Range("A81").Value = -22.24240000000001
$ws.Range("A82").Value = -21.7011
$ws.Range("A84").Value = -21.9947
$ws.Range("B84").Value = 5.9479
$ws.Range("B85").Value = 5.9654
$ws.Range("B89").Value = 4.351699999999997
$ws.Range("B93").Value = 5.451599999999997
$ws.Range("B95").Value = 5.401000000000002
$ws.Range("B98").Value = 7.705000000000001
$ws.Range("B99").Value = 5.820899999999997
$ws.Range("A101").Value = -21.72449999999999
$ws.Range("B101").Value = 5.388799999999994
$ws.Range("A102").Value = -22.0667
$ws.Range("B102").Value = 5.9518
